# Added logistics fleet as pure consumption asset, integrated evaluation in
# Building gridConnection.

$wb = $excel.ActiveWorkbook

# --- consumptionAssets: append new row for the logistics fleet HGV asset ---
$wsCons = $wb.Worksheets.Item("consumptionAssets")

$wsCons.Range("A10").Value = 9
$wsCons.Range("B10").Value = "Logistics_fleet_hgv_E"
$wsCons.Range("C10").Value = "logistics_fleet_e_hgv"
$wsCons.Range("D10").Value = "CONSUMPTION"
$wsCons.Range("E10").Value = "ELECTRICITY_DEMAND"
$wsCons.Range("F10").Formula = "=25*100000*3"
$wsCons.Range("G10").Value = 0

# --- storageAssets: rename the generic grid battery to a sized 1 MWh variant ---
$wsStor = $wb.Worksheets.Item("storageAssets")
$wsStor.Range("B14").Value = "Grid_battery_1MWh"

# --- selections / active sheet, matching the editor's final on-screen state ---
[void]$wsStor.Range("B14").Select()

[void]$wsCons.Activate()
[void]$wsCons.Range("F10").Select()
